# Fix logging system configuration
# Appends a new log row (row 75) to each of the four sheets, mirroring the
# prior row (row 74) but advancing the timestamp in column A by ~1 day.

$wb = $excel.ActiveWorkbook

$rowsData = @(
    @{ Sheet = "FE_LFT_#1"; A = 45861.49177083333; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x2C"; E = "0xf"; F = 380; G = "7.598631275147109e+23"; H = 300; I = 15 }
    @{ Sheet = "FE_LFT_#2"; A = 45861.49177083333; B = "0x01,0x90"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x3C"; E = "0xe"; F = 400; G = "5.68432987514711e+23";  H = 316; I = 14 }
    @{ Sheet = "FE_PLT_#1"; A = 45861.49177083333; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x61"; E = "0x3"; F = 110; G = "5.68631262647114e+23"; H = 97;  I = 3 }
    @{ Sheet = "FE_PLT_#2"; A = 45861.49177083333; B = "0x00,0x6e"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x61"; E = "0x3"; F = 110; G = "9.85046333984776e+23"; H = 97;  I = 3 }
)

foreach ($rd in $rowsData) {
    $ws = $wb.Worksheets.Item($rd.Sheet)
    $newRow = 75

    $ws.Cells.Item($newRow, 1).Value = $rd.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $rd.B
    $ws.Cells.Item($newRow, 3).Value = $rd.C
    $ws.Cells.Item($newRow, 4).Value = $rd.D
    $ws.Cells.Item($newRow, 5).Value = $rd.E

    $ws.Cells.Item($newRow, 6).Value = $rd.F
    $ws.Cells.Item($newRow, 7).Value = [double]$rd.G
    $ws.Cells.Item($newRow, 8).Value = $rd.H
    $ws.Cells.Item($newRow, 9).Value = $rd.I
}
